# Apply NATMI Pdgfb-Pdgfra update (following Dr Hou advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2-4: refresh A/D labels and numeric columns for existing rows, extend B/C ---
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Pdgfb"
$ws.Cells.Item(2, 3).Value = "Pdgfra"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 41.23709466666667
$ws.Cells.Item(2, 8).Value = 123.711284
$ws.Cells.Item(2, 9).Value = 0.956365997213294
$ws.Cells.Item(2, 10).Value = 0.9563659972132939
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.2781686666666667
$ws.Cells.Item(2, 14).Value = 0.834506
$ws.Cells.Item(2, 15).Value = 0.001228014730390642
$ws.Cells.Item(2, 16).Value = 0.001228014730390642
$ws.Cells.Item(2, 17).Value = 11.47086764063378
$ws.Cells.Item(2, 18).Value = 103.237808765704
$ws.Cells.Item(2, 19).Value = 0.001174431532222661
$ws.Cells.Item(2, 20).Value = 0.001174431532222661

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Pdgfb"
$ws.Cells.Item(3, 3).Value = "Pdgfra"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 41.23709466666667
$ws.Cells.Item(3, 8).Value = 123.711284
$ws.Cells.Item(3, 9).Value = 0.956365997213294
$ws.Cells.Item(3, 10).Value = 0.9563659972132939
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 225.778076
$ws.Cells.Item(3, 14).Value = 677.3342279999999
$ws.Cells.Item(3, 15).Value = 0.9967290940769435
$ws.Cells.Item(3, 16).Value = 0.9967290940769435
$ws.Cells.Item(3, 17).Value = 9310.431893669862
$ws.Cells.Item(3, 18).Value = 83793.88704302875
$ws.Cells.Item(3, 19).Value = 0.9532378140083992
$ws.Cells.Item(3, 20).Value = 0.9532378140083991

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Pdgfb"
$ws.Cells.Item(4, 3).Value = "Pdgfra"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 41.23709466666667
$ws.Cells.Item(4, 8).Value = 123.711284
$ws.Cells.Item(4, 9).Value = 0.956365997213294
$ws.Cells.Item(4, 10).Value = 0.9563659972132939
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.4627536666666667
$ws.Cells.Item(4, 14).Value = 1.388261
$ws.Cells.Item(4, 15).Value = 0.002042891192665893
$ws.Cells.Item(4, 16).Value = 0.002042891192665893
$ws.Cells.Item(4, 17).Value = 19.08261675968045
$ws.Cells.Item(4, 18).Value = 171.743550837124
$ws.Cells.Item(4, 19).Value = 0.001953751672672172
$ws.Cells.Item(4, 20).Value = 0.001953751672672172

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Pdgfb"
$ws.Cells.Item(5, 3).Value = "Pdgfra"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.030543
$ws.Cells.Item(5, 8).Value = 0.091629
$ws.Cells.Item(5, 9).Value = 0.0007083497731593903
$ws.Cells.Item(5, 10).Value = 0.0007083497731593902
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.2781686666666667
$ws.Cells.Item(5, 14).Value = 0.834506
$ws.Cells.Item(5, 15).Value = 0.001228014730390642
$ws.Cells.Item(5, 16).Value = 0.001228014730390642
$ws.Cells.Item(5, 17).Value = 0.008496105586
$ws.Cells.Item(5, 18).Value = 0.076464950274
$ws.Cells.Item(5, 19).Value = 0.0000008698639557086014
$ws.Cells.Item(5, 20).Value = 0.0000008698639557086012

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Pdgfb"
$ws.Cells.Item(6, 3).Value = "Pdgfra"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.030543
$ws.Cells.Item(6, 8).Value = 0.091629
$ws.Cells.Item(6, 9).Value = 0.0007083497731593903
$ws.Cells.Item(6, 10).Value = 0.0007083497731593902
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 225.778076
$ws.Cells.Item(6, 14).Value = 677.3342279999999
$ws.Cells.Item(6, 15).Value = 0.9967290940769435
$ws.Cells.Item(6, 16).Value = 0.9967290940769435
$ws.Cells.Item(6, 17).Value = 6.895939775267999
$ws.Cells.Item(6, 18).Value = 62.06345797741199
$ws.Cells.Item(6, 19).Value = 0.0007060328276907676
$ws.Cells.Item(6, 20).Value = 0.0007060328276907675

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Pdgfb"
$ws.Cells.Item(7, 3).Value = "Pdgfra"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.030543
$ws.Cells.Item(7, 8).Value = 0.091629
$ws.Cells.Item(7, 9).Value = 0.0007083497731593903
$ws.Cells.Item(7, 10).Value = 0.0007083497731593902
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.4627536666666667
$ws.Cells.Item(7, 14).Value = 1.388261
$ws.Cells.Item(7, 15).Value = 0.002042891192665893
$ws.Cells.Item(7, 16).Value = 0.002042891192665893
$ws.Cells.Item(7, 17).Value = 0.014133885241
$ws.Cells.Item(7, 18).Value = 0.127204967169
$ws.Cells.Item(7, 19).Value = 0.000001447081512914202
$ws.Cells.Item(7, 20).Value = 0.000001447081512914201

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Pdgfb"
$ws.Cells.Item(8, 3).Value = "Pdgfra"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.850891
$ws.Cells.Item(8, 8).Value = 5.552673
$ws.Cells.Item(8, 9).Value = 0.0429256530135467
$ws.Cells.Item(8, 10).Value = 0.04292565301354669
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.2781686666666667
$ws.Cells.Item(8, 14).Value = 0.834506
$ws.Cells.Item(8, 15).Value = 0.001228014730390642
$ws.Cells.Item(8, 16).Value = 0.001228014730390642
$ws.Cells.Item(8, 17).Value = 0.5148598816153334
$ws.Cells.Item(8, 18).Value = 4.633738934537999
$ws.Cells.Item(8, 19).Value = 0.00005271333421227281
$ws.Cells.Item(8, 20).Value = 0.00005271333421227281

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Pdgfb"
$ws.Cells.Item(9, 3).Value = "Pdgfra"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.850891
$ws.Cells.Item(9, 8).Value = 5.552673
$ws.Cells.Item(9, 9).Value = 0.0429256530135467
$ws.Cells.Item(9, 10).Value = 0.04292565301354669
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 225.778076
$ws.Cells.Item(9, 14).Value = 677.3342279999999
$ws.Cells.Item(9, 15).Value = 0.9967290940769435
$ws.Cells.Item(9, 16).Value = 0.9967290940769435
$ws.Cells.Item(9, 17).Value = 417.8906088657159
$ws.Cells.Item(9, 18).Value = 3761.015479791443
$ws.Cells.Item(9, 19).Value = 0.04278524724085362
$ws.Cells.Item(9, 20).Value = 0.04278524724085361

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Pdgfb"
$ws.Cells.Item(10, 3).Value = "Pdgfra"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.850891
$ws.Cells.Item(10, 8).Value = 5.552673
$ws.Cells.Item(10, 9).Value = 0.0429256530135467
$ws.Cells.Item(10, 10).Value = 0.04292565301354669
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.4627536666666667
$ws.Cells.Item(10, 14).Value = 1.388261
$ws.Cells.Item(10, 15).Value = 0.002042891192665893
$ws.Cells.Item(10, 16).Value = 0.002042891192665893
$ws.Cells.Item(10, 17).Value = 0.8565065968503333
$ws.Cells.Item(10, 18).Value = 7.708559371652999
$ws.Cells.Item(10, 19).Value = 0.0000876924384808067
$ws.Cells.Item(10, 20).Value = 0.00008769243848080666
